$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(1, "Subsector", "Comment"),
    @(2, "All Subsectors", "1 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(3, "All Transport Subsectors", "2 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(4, "All Industry Subsectors", "3 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(5, "All Commercial Subsectors", "4 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(6, "All Residential Subsectors", "5 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(7, "All Agriculture Subsectors", "6 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(8, "All Other Subsectors", "7 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(9, "Aluminium", "2 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(10, "Aviation", "Emissions from international aviation are not constraint in the model. Fuel consumption for aviation decreases in the period 2045-2060 in the Kea scenario because we assumed an aspect of flight shame."),
    @(11, "Chemical", "4 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(12, "Construction", "5 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(13, "Dairy", "6 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(14, "Dairy Cattle Farming", "7 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(15, "Detached Dwellings", "Solid fuels (wood and coal) are only considered in detached dwellings. Natural Gas and Geothermal energy are only considered in the North Island."),
    @(16, "Distributed Battery", "9 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(17, "Education", "10 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(18, "Fishing", "11 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(19, "Food (Non Dairy/Meat)", "12 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(20, "Forestry", "13 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(21, "Geothermal", "14 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(22, "Healthcare", "15 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(23, "Heavy Road", "Different learning curves, and therefore costs, are assumed for electric and hydrogen technologies in each scenario."),
    @(24, "Horticulture", "17 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(25, "Hydro", "18 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(26, "Indoor Cropping", "19 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(27, "Iron/Steel", "20 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(28, "Joint Dwellings", "Joint dwellings comprises medium-density dwellings and apartments. Natural Gas and Geothermal energy are only considered in the North Island."),
    @(29, "Light Road", "In each time period, the share of EVs is larger in Kea because we assume that in such scenario there is a larger ability to access these cars."),
    @(30, "Livestock Farming", "23 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(31, "Meat", "24 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(32, "Metal", "25 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(33, "Methanol", "26 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(34, "Minerals", "27 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(35, "Mining", "28 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(36, "Office Blocks", "29 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(37, "Other", "30 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(38, "Pulp and paper", "31 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(39, "Pumped Hydro", "32 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(40, "Rail", "Passenger rail in the South Island was disregarded because it represents only a minor amount."),
    @(41, "Refining", "34 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(42, "Shipping", "Emissions from international shipping are not constraint in the model."),
    @(43, "Solar", "36 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(44, "Thermal", "37 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(45, "Warehouses/Supermarkets/Retail", "38 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(46, "Wind", "39 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  "),
    @(47, "Wood", "40 Draft commentary: The scenario shows that petrol has high consumption until 2035 at which point in sharply decreases  ")

)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

$ws.Columns.Item(1).ColumnWidth = 32.666666666666664
$ws.Columns.Item(2).ColumnWidth = 69.83333333333333

$ws.Range("E5").Select()
